$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header label for the strata column, and rename the other two headers
$ws.Range("A1").Value = "Strata"
$ws.Range("B1").Value = "Area"
$ws.Range("C1").Value = "SE"

# Move active selection to D1 (reflects user's final click position)
$ws.Range("D1").Select()
